$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the "District of Columbia" row (row 12) -- all rows below shift up.
$ws.Rows.Item(12).Delete()
